$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 249.2
$ws.Range("I6").Value = 23
$ws.Range("K6").Value = 69
$ws.Range("M6").Value = 43
$ws.Range("H9").Value = 84.833336
$ws.Range("I9").Value = 76
$ws.Range("J9").Value = 129
$ws.Range("K9").Value = 76
$ws.Range("L9").Value = 129
$ws.Range("M9").Value = 93
$ws.Range("N9").Value = -467
$ws.Range("H13").Value = 2011
$ws.Range("I13").Value = 1014.6667
$ws.Range("K13").Value = 1014.6667
$ws.Range("M13").Value = -845.6667
$ws.Range("H15").Value = 784
$ws.Range("I15").Value = 784
$ws.Range("K15").Value = 2352
$ws.Range("M15").Value = -2183
$ws.Range("H16").Value = 550
$ws.Range("I16").Value = 550
$ws.Range("K16").Value = 550
$ws.Range("M16").Value = -320
$ws.Range("H18").Value = 1969.3334
$ws.Range("J18").Value = 2454
$ws.Range("L18").Value = 2454
$ws.Range("N18").Value = -3022
$ws.Range("H33").Value = 264.9524
$ws.Range("I33").Value = 272.7
$ws.Range("J33").Value = 110
$ws.Range("K33").Value = 272.7
$ws.Range("L33").Value = 110
$ws.Range("M33").Value = -43.69999999999999
$ws.Range("N33").Value = -568
$ws.Range("H34").Value = 14999.5
$ws.Range("I34").Value = 14999.5
$ws.Range("K34").Value = 14999.5
$ws.Range("M34").Value = -14796.5
$ws.Range("H36").Value = 14999.5
$ws.Range("I36").Value = 14999.5
$ws.Range("K36").Value = 14999.5
$ws.Range("M36").Value = -14284.5
$ws.Range("H38").Value = 71.77778000000001
$ws.Range("I38").Value = 71.77778000000001
$ws.Range("K38").Value = 215.33334
$ws.Range("M38").Value = 156.66666
$ws.Range("H62").Value = 8149.6665
$ws.Range("I62").Value = 2199.5
$ws.Range("J62").Value = 11124.75
$ws.Range("K62").Value = 2199.5
$ws.Range("L62").Value = 11124.75
$ws.Range("M62").Value = -1575.5
$ws.Range("N62").Value = -12372.75
$ws.Range("H64").Value = 4500
$ws.Range("J64").Value = 4500
$ws.Range("L64").Value = 4500
$ws.Range("N64").Value = -4996
$ws.Range("H65").Value = 8149.6665
$ws.Range("I65").Value = 2199.5
$ws.Range("J65").Value = 11124.75
$ws.Range("K65").Value = 10997.5
$ws.Range("L65").Value = 55623.75
$ws.Range("M65").Value = -7877.5
$ws.Range("N65").Value = -61863.75
$ws.Range("H67").Value = 4500
$ws.Range("J67").Value = 4500
$ws.Range("L67").Value = 4500
$ws.Range("N67").Value = -6216
$ws.Range("H98").Value = 5085
$ws.Range("I98").Value = 3152.5
$ws.Range("J98").Value = 8950
$ws.Range("K98").Value = 3152.5
$ws.Range("L98").Value = 8950
$ws.Range("M98").Value = -1654.5
$ws.Range("N98").Value = -11946
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 5085
$ws.Range("I122").Value = 3152.5
$ws.Range("J122").Value = 8950
$ws.Range("K122").Value = 9457.5
$ws.Range("L122").Value = 26850
$ws.Range("M122").Value = -7007.5
$ws.Range("N122").Value = -31750

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H117").Value = 180000
$ws.Range("J117").Value = 180000
$ws.Range("L117").Value = 180000
$ws.Range("N117").Value = -189178
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 6606.25
$ws.Range("I132").Value = 2262.5
$ws.Range("K132").Value = 6787.5
$ws.Range("M132").Value = -4257.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 39890
$ws.Range("J81").Value = 39890
$ws.Range("L81").Value = 39890
$ws.Range("N81").Value = -42012
$ws.Range("H84").Value = 39890
$ws.Range("J84").Value = 39890
$ws.Range("L84").Value = 119670
$ws.Range("N84").Value = -130278
$ws.Range("H94").Value = 1828.1666
$ws.Range("I94").Value = 1828.1666
$ws.Range("K94").Value = 1828.1666
$ws.Range("M94").Value = -1377.1666
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("H140").Value = 49997
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 49997
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 49997
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -60357

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3820.8333
$ws.Range("J31").Value = 4733.75
$ws.Range("L31").Value = 4733.75
$ws.Range("N31").Value = -5323.75
$ws.Range("H34").Value = 3820.8333
$ws.Range("J34").Value = 4733.75
$ws.Range("L34").Value = 4733.75
$ws.Range("N34").Value = -5137.75
$ws.Range("H41").Value = 1250
$ws.Range("I41").Value = 1250
$ws.Range("K41").Value = 1250
$ws.Range("M41").Value = -822
$ws.Range("H58").Value = 3966.3333
$ws.Range("I58").Value = 1900
$ws.Range("J58").Value = 4999.5
$ws.Range("K58").Value = 1900
$ws.Range("L58").Value = 4999.5
$ws.Range("M58").Value = -1697
$ws.Range("N58").Value = -5405.5
$ws.Range("H63").Value = 30000
$ws.Range("J63").Value = 30000
$ws.Range("L63").Value = 30000
$ws.Range("N63").Value = -31372
$ws.Range("H66").Value = 30000
$ws.Range("J66").Value = 30000
$ws.Range("L66").Value = 90000
$ws.Range("N66").Value = -96864
$ws.Range("H132").Value = 5402.4
$ws.Range("I132").Value = 5402.4
$ws.Range("K132").Value = 16207.2
$ws.Range("M132").Value = -13677.2
$ws.Range("H136").Value = 3966.3333
$ws.Range("I136").Value = 1900
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 5700
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = -3150
$ws.Range("N136").Value = -20098.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H34").Value = 1928.75
$ws.Range("I34").Value = 350
$ws.Range("J34").Value = 2072.2727
$ws.Range("K34").Value = 1050
$ws.Range("L34").Value = 6216.8181
$ws.Range("M34").Value = -966
$ws.Range("N34").Value = -6384.8181
$ws.Range("H97").Value = 351.66666
$ws.Range("I97").Value = 152.5
$ws.Range("J97").Value = 750
$ws.Range("K97").Value = 457.5
$ws.Range("L97").Value = 2250
$ws.Range("M97").Value = 38.5
$ws.Range("N97").Value = -3242
$ws.Range("H98").Value = 651.3333
$ws.Range("I98").Value = 699.5
$ws.Range("J98").Value = 555
$ws.Range("K98").Value = 2098.5
$ws.Range("L98").Value = 1665
$ws.Range("M98").Value = -600.5
$ws.Range("N98").Value = -4661
$ws.Range("H112").Value = 1975
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H113").Value = 978.125
$ws.Range("I113").Value = 644.25
$ws.Range("J113").Value = 1312
$ws.Range("K113").Value = 1932.75
$ws.Range("L113").Value = 3936
$ws.Range("M113").Value = 237.25
$ws.Range("N113").Value = -8276
$ws.Range("H122").Value = 100
$ws.Range("I122").Value = 100
$ws.Range("J122").Value = 100
$ws.Range("K122").Value = 900
$ws.Range("L122").Value = 900
$ws.Range("M122").Value = 1550
$ws.Range("N122").Value = -5800

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 110
$ws.Range("I9").Value = 110
$ws.Range("K9").Value = 110
$ws.Range("M9").Value = 60
$ws.Range("H101").Value = 49999
$ws.Range("J101").Value = 49999
$ws.Range("L101").Value = 49999
$ws.Range("N101").Value = -56489
$ws.Range("H113").Value = 1111
$ws.Range("I113").Value = 1111
$ws.Range("K113").Value = 1111
$ws.Range("M113").Value = 1059
$ws.Range("H132").Value = 2847.25
$ws.Range("I132").Value = 2924.3635
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 8773.0905
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -6243.0905
$ws.Range("N132").Value = -11057
$ws.Range("H140").Value = 47139.715
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 47139.715
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 47139.715
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -57499.715

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 168.5
$ws.Range("I10").Value = 168.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 168.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -28.5
$ws.Range("N10").ClearContents()
$ws.Range("H22").Value = 1049.5
$ws.Range("I22").Value = 732.6667
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 732.6667
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -437.6667
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1049.5
$ws.Range("I27").Value = 732.6667
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 732.6667
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -625.6667
$ws.Range("N27").Value = -2214
$ws.Range("H41").Value = 6000
$ws.Range("I41").Value = 6000
$ws.Range("K41").Value = 6000
$ws.Range("M41").Value = -5562

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5000
$ws.Range("J107").Value = 5000
$ws.Range("L107").Value = 15000
$ws.Range("N107").Value = -18840
$ws.Range("H122").Value = 2631.6667
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 999
$ws.Range("I132").Value = 999
$ws.Range("K132").Value = 2997
$ws.Range("M132").Value = -467

